# Applies the weekly refresh of Fruta/Hortalizas price data for the
# "Terminal Hortofrutícola Agro Chillán - Mango" subset.
# The full 40-row dataset is re-sorted (dates shuffled) and 3 new price
# records are appended (rows 42-44), extending the sheet to A1:T44.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that stay constant across every data row.
$mercadoId = 7
$mercado   = "Terminal Hortofrutícola Agro Chillán"
$region    = "Ñuble"
$codreg    = 16
$tipo      = "Fruta"
$productoId = 100108
$producto  = "Tropicales y subtropicales"
$categoriaId = 100108002
$categoria = "Mango"
$variedad  = "Sin especificar"
$unidad    = "$/bandeja 4 kilos"
$kgUnidad  = 4

# Per-row data: Fecha(serial), Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, Origen, PrecioKg
$rows = @(
    @(44211, "Primera", 50, 6000, 6000, 6000, "Perú", 1500),
    @(44424, "Primera", 50, 8500, 9000, 8800, "Perú", 2200),
    @(44214, "Especial", 40, 6000, 6000, 6000, "Perú", 1500),
    @(44214, "Primera", 50, 6000, 6000, 6000, "Perú", 1500),
    @(44188, "Primera", 50, 6000, 6000, 6000, "Perú", 1500),
    @(44445, "Primera", 60, 8500, 9000, 8750, "Brasil", 2188),
    @(44232, "Primera", 60, 6000, 6000, 6000, "Perú", 1500),
    @(44216, "Primera", 50, 6000, 6000, 6000, "Perú", 1500),
    @(44229, "Especial", 25, 6000, 6000, 6000, "Perú", 1500),
    @(44229, "Primera", 35, 6000, 6000, 6000, "Perú", 1500),
    @(44189, "Primera", 40, 6000, 6000, 6000, "Perú", 1500),
    @(44181, "Primera", 50, 6000, 6000, 6000, "Perú", 1500),
    @(44181, "Segunda", 40, 5000, 5000, 5000, "Perú", 1250),
    @(44231, "Especial", 45, 6500, 6500, 6500, "Perú", 1625),
    @(44231, "Primera", 35, 6500, 6500, 6500, "Perú", 1625),
    @(44230, "Especial", 45, 6500, 6500, 6500, "Perú", 1625),
    @(44230, "Primera", 50, 6000, 6000, 6000, "Perú", 1500),
    @(44174, "Especial", 50, 6500, 7000, 6800, "Perú", 1700),
    @(44174, "Primera", 70, 6000, 7000, 6286, "Perú", 1572),
    @(44187, "Primera", 40, 6000, 6000, 6000, "Perú", 1500),
    @(44195, "Primera", 50, 6000, 6000, 6000, "Perú", 1500),
    @(44209, "Primera", 50, 6000, 6000, 6000, "Perú", 1500),
    @(44186, "Primera", 40, 6000, 6000, 6000, "Perú", 1500),
    @(44252, "Primera", 50, 6000, 6000, 6000, "Perú", 1500),
    @(44419, "Primera", 60, 8500, 9000, 8750, "Perú", 2188),
    @(44292, "Primera", 30, 7500, 8000, 7750, "Perú", 1938),
    @(44179, "Primera", 50, 7000, 7000, 7000, "Perú", 1750),
    @(44179, "Segunda", 40, 6000, 6000, 6000, "Perú", 1500),
    @(44196, "Primera", 50, 6000, 6000, 6000, "Perú", 1500),
    @(44218, "Primera", 40, 6000, 6000, 6000, "Perú", 1500),
    @(44433, "Primera", 60, 8500, 9000, 8750, "Perú", 2188),
    @(44426, "Primera", 40, 8500, 9000, 8750, "Perú", 2188),
    @(44222, "Primera", 50, 6000, 6000, 6000, "Perú", 1500),
    @(44420, "Primera", 40, 8500, 9000, 8750, "Perú", 2188),
    @(44210, "Primera", 40, 6000, 6000, 6000, "Perú", 1500),
    @(44253, "Especial", 30, 6000, 6000, 6000, "Perú", 1500),
    @(44253, "Primera", 50, 6000, 6000, 6000, "Perú", 1500),
    @(44417, "Primera", 80, 8500, 9000, 8750, "Perú", 2188),
    @(44446, "Primera", 60, 8500, 9000, 8750, "Brasil", 2188),
    @(44438, "Primera", 60, 8500, 9000, 8750, "Perú", 2188),
    @(44251, "Primera", 65, 6000, 6500, 6269, "Perú", 1567),
    @(44442, "Primera", 60, 8500, 9000, 8750, "Ecuador", 2188),
    @(44435, "Primera", 60, 8500, 9000, 8750, "Perú", 2188)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region

    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Value = $data[0]
    $dCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $productoId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $data[1]
    $ws.Cells.Item($r, 13).Value = $data[2]
    $ws.Cells.Item($r, 14).Value = $data[3]
    $ws.Cells.Item($r, 15).Value = $data[4]
    $ws.Cells.Item($r, 16).Value = $data[5]
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $data[6]
    $ws.Cells.Item($r, 19).Value = $data[7]
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}
